$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New prediction-history rows being appended below the existing data (rows 2-7).
# Amount and Risk Score must remain plain text (e.g. "$9,450.00", "0.08%") exactly
# like the rest of the sheet, instead of being auto-parsed into currency/percent
# numbers, so each of those two cells is briefly marked as Text ("@") before the
# value is assigned and then has ClearFormats() applied so no residual number
# format / style lingers on the cell.
$rows = @(
    @{ Row=8;  Ts="2025-11-25 02:45:24.994859"; Sender=800072000; Receiver=100428660; Amount="`$9,450.00";         Currency="UK Pound"; Format="ACH"; Risk="0.08%"; Status="LOW RISK" },
    @{ Row=9;  Ts="2025-11-25 02:46:55.004932"; Sender=800072000; Receiver=100428660; Amount="`$9,450.00";         Currency="UK Pound"; Format="ACH"; Risk="0.08%"; Status="LOW RISK" },
    @{ Row=10; Ts="2025-11-25 02:47:34.308052"; Sender=800072000; Receiver=100428660; Amount="`$9,450.00";         Currency="UK Pound"; Format="ACH"; Risk="0.08%"; Status="LOW RISK" },
    @{ Row=11; Ts="2025-11-22 02:49:59.295945"; Sender=800072000; Receiver=100428660; Amount="`$150,000.00";       Currency="UK Pound"; Format="ACH"; Risk="0.08%"; Status="LOW RISK" },
    @{ Row=12; Ts="2025-11-15 02:53:52.912504"; Sender=800072000; Receiver=100428660; Amount="`$1,000,000.00";     Currency="UK Pound"; Format="ACH"; Risk="0.08%"; Status="LOW RISK" },
    @{ Row=13; Ts="2025-11-24 02:56:55.767512"; Sender=1004;      Receiver=2156;      Amount="`$9,500.00";         Currency="UK Pound"; Format="ACH"; Risk="0.23%"; Status="LOW RISK" },
    @{ Row=14; Ts="2025-11-24 02:57:05.429107"; Sender=1004;      Receiver=2156;      Amount="`$9,500.00";         Currency="UK Pound"; Format="ACH"; Risk="0.23%"; Status="LOW RISK" },
    @{ Row=15; Ts="2022-11-24 02:58:56.682773"; Sender=1004;      Receiver=2156;      Amount="`$9,500.00";         Currency="UK Pound"; Format="ACH"; Risk="0.13%"; Status="LOW RISK" },
    @{ Row=16; Ts="2022-11-25 03:00:51.456492"; Sender=1004;      Receiver=2156;      Amount="`$9,500.00";         Currency="UK Pound"; Format="ACH"; Risk="0.19%"; Status="LOW RISK" },
    @{ Row=17; Ts="2025-11-25 03:08:44.898061"; Sender=800072000; Receiver=100428660; Amount="`$9,450.00";         Currency="UK Pound"; Format="ACH"; Risk="0.08%"; Status="LOW RISK" }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("A$n").Value = $r.Ts
    $ws.Range("B$n").Value = $r.Sender
    $ws.Range("C$n").Value = $r.Receiver

    $ws.Range("D$n").NumberFormat = "@"
    $ws.Range("D$n").Value = $r.Amount
    $ws.Range("D$n").ClearFormats()

    $ws.Range("E$n").Value = $r.Currency
    $ws.Range("F$n").Value = $r.Format

    $ws.Range("G$n").NumberFormat = "@"
    $ws.Range("G$n").Value = $r.Risk
    $ws.Range("G$n").ClearFormats()

    $ws.Range("H$n").Value = $r.Status
}
